# Apply the "Trade #33 closed at 2026-02-17 20:57:36" update to the
# live trading results workbook.
#
# Sheets:
#   Summary          -> sheet1
#   Strategy Status   -> sheet2
#   All Trades        -> sheet3
#   MarketMaking      -> sheet4

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet - update aggregate stats
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.46   # Current Capital
$summary.Range("B4").Value = 0.25      # Total P&L $
$summary.Range("B5").Value = 0.08      # Total P&L %
$summary.Range("B6").Value = 61        # Total Trades
$summary.Range("B8").Value = 24        # Losing Trades
$summary.Range("B9").Value = 47.54     # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - update MarketMaking row (row 5)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.46     # Capital
$status.Range("D5").Value = 28         # Trades
$status.Range("E5").Value = 0.14       # P&L $
$status.Range("F5").Value = 0.46       # P&L %
$status.Range("G5").Value = 53.57      # Win Rate %

# ---------------------------------------------------------------------
# 3) All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Trade #61 (row 62) transitions from OPEN to CLOSED via early_exit
$allTrades.Range("G62").Value = 0.8
$allTrades.Range("H62").Value = "CLOSED"
$allTrades.Range("I62").Value = -4.7619
$allTrades.Range("J62").Value = -0.04
$allTrades.Range("K62").Value = 100.46
$allTrades.Range("L62").Value = "early_exit"
$allTrades.Range("M62").Value = 0.14

# New trade (Trade #94) opened, appended as row 95
$allTrades.Range("A95").Value = 94
$allTrades.Range("B95").NumberFormat = "@"
$allTrades.Range("B95").Value = "2026-02-17"
$allTrades.Range("C95").Value = "20:57:29"
$allTrades.Range("D95").Value = "MarketMaking"
$allTrades.Range("E95").Value = "DOWN"
$allTrades.Range("F95").Value = 0.84
$allTrades.Range("H95").Value = "OPEN"
$allTrades.Range("I95").Value = 0
$allTrades.Range("J95").Value = 0
$allTrades.Range("K95").Value = 100.4955022889912
$allTrades.Range("M95").Value = 0
$allTrades.Range("N95").Value = 0
$allTrades.Range("O95").Value = 0
$allTrades.Range("P95").Value = 0.6
$allTrades.Range("Q95").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# 4) MarketMaking sheet (strategy-specific trade log)
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# Same Trade #61 (row 29) transitions from OPEN to CLOSED via early_exit
$mm.Range("G29").Value = 0.8
$mm.Range("H29").Value = "CLOSED"
$mm.Range("I29").Value = -4.7619
$mm.Range("J29").Value = -0.04
$mm.Range("K29").Value = 100.46
$mm.Range("P29").Value = "early_exit"
$mm.Range("Q29").Value = 0.14

# New trade (Trade #94) opened, appended as row 62
$mm.Range("A62").Value = 94
$mm.Range("B62").NumberFormat = "@"
$mm.Range("B62").Value = "2026-02-17"
$mm.Range("C62").Value = "20:57:29"
$mm.Range("D62").Value = "MarketMaking"
$mm.Range("E62").Value = "DOWN"
$mm.Range("F62").Value = 0.84
$mm.Range("H62").Value = "OPEN"
$mm.Range("I62").Value = 0
$mm.Range("J62").Value = 0
$mm.Range("K62").Value = 100.4955022889912
$mm.Range("L62").Value = 0
$mm.Range("M62").Value = 0
$mm.Range("N62").Value = 0.6
$mm.Range("O62").Value = "Normal spread capture: 19600 bps"
$mm.Range("Q62").Value = 0

Write-Output "Applied trade #33/61 close + new trade #94 open across Summary, Strategy Status, All Trades, MarketMaking sheets."
